$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated, in order: B, C, D, E, F, H, I, J, K, M, N
$cols = @(2,3,4,5,6,8,9,10,11,13,14)

# New values for rows 2..25 (24 rows), one row per sub-array, values in $cols order
$data = @(
    @(13.62872295045841,10.69096673853173,6.276190333643999,12.04723991820001,45.50761799847226,7.344005520526261,30.69713944760969,10.22790188619332,14.81760398917559,17.71410498327277,23.10245998556931),
    @(13.45722975573768,10.57631360748534,6.269605433283242,12.05130322508267,45.48662517687423,7.344005520526261,30.73301094569063,10.24619272929466,14.71138967262527,17.68454041789051,23.15580052701283),
    @(13.35452903539435,10.50828343281539,6.266529836636754,12.05570197436574,45.48284905741058,7.344005520526261,30.76015662726954,10.25852432375589,14.64942813518384,17.6698002950344,23.19046458070335),
    @(13.31338244821324,10.48118672906387,6.265521279454783,12.05797392139187,45.48360349104158,7.344005520526261,30.77250464953034,10.26382665022754,14.62502032929406,17.66465637923096,23.20507201218714),
    @(13.3065941424708,10.47672603467379,6.265368633587262,12.05838014925637,45.48386728991926,7.344005520526261,30.77463264334821,10.26472384128545,14.62101894965718,17.66385447089107,23.20752666645696),
    @(13.3539711970769,10.50791542220347,6.266515241942443,12.05573067268878,45.48284994575426,7.344005520526261,30.76031795274289,10.25859471047214,14.64909552353856,17.66972742325439,23.19065963113122),
    @(13.56908461594819,10.65096192424879,6.273719962790931,12.04824624827998,45.4984894546978,7.344005520526261,30.70844428241123,10.23398022626683,14.78032046734931,17.70320581625936,23.12045515821158),
    @(14.00919210843043,10.94884274186859,6.295456120890998,12.04864105365547,45.60135339612515,7.344005520526261,30.64741379164501,10.19443533245918,15.06236727845422,17.79570304893639,22.99793482200434),
    @(14.3405038778883,11.17631988751178,6.315963619154449,12.0580658506187,45.72071058970721,7.344005520526261,30.62745983616417,10.1706846859703,15.2830834239657,17.87967844615647,22.91711806616622),
    @(14.49227977336341,11.28124851244456,6.326254530029013,12.06432241226988,45.7844399811915,7.344005520526261,30.62379643886287,10.16102802841642,15.38606779861141,17.92127179712254,22.8823422746757),
    @(14.54985007968512,11.32115428197604,6.330287523231616,12.06697323346824,45.80992020810182,7.344005520526261,30.62318792624907,10.15753604771355,15.42540521260341,17.93750154121341,22.86945898187838),
    @(14.53744789852866,11.31255281130741,6.329412932750612,12.0663898343351,45.80437282529643,7.344005520526261,30.6232843455089,10.15828078327838,15.41691860505752,17.93398500407339,22.87222093515813),
    @(14.49701451119293,11.28452835481745,6.326583613668965,12.06453486287503,45.78650931066914,7.344005520526261,30.62373076893096,10.16073743992181,15.3892974899288,17.92259747450528,22.88127663761495),
    @(14.4722587480998,11.26738382522032,6.324868221899699,12.06343526216071,45.77574255261967,7.344005520526261,30.62410563211801,10.16226366572373,15.37242202540736,17.91568443038052,22.88686069316153),
    @(14.33060145947367,11.1694885579454,6.315310248793607,12.05769646423157,45.71673478991264,7.344005520526261,30.62780821505675,10.17133884401649,15.27640238179135,17.87702776796262,22.91943072234554),
    @(14.24393046568118,11.1097777620956,6.309691537446047,12.05467918810088,45.68294540673269,7.344005520526261,30.6314665040787,10.17719993810478,15.21813590415548,17.85417613926873,22.9399203395446),
    @(14.19418274715003,11.07557218989582,6.306550490491691,12.05312918313306,45.66439918325086,7.344005520526261,30.63408023001751,10.18067911979567,15.18486728476961,17.84135236021178,22.95189257957915),
    @(14.17735847525451,11.06401567726396,6.305502628986669,12.05263627474195,45.65827261516257,7.344005520526261,30.63505269977362,10.18187567336286,15.17364608193938,17.83706563960059,22.9559783292834),
    @(14.25314647783613,11.11612001024475,6.310280289201126,12.05498119949226,45.68645044624289,7.344005520526261,30.63102433097561,10.17656483510016,15.22431336902183,17.8565756888155,22.93771981762217),
    @(14.50888861433496,11.29275545174468,6.327410979702903,12.06507208381893,45.79171977297972,7.344005520526261,30.62357850859263,10.16001139004107,15.39740152120536,17.92592933286066,22.87860901137243),
    @(14.67656549106376,11.40918017596743,6.33939873035468,12.07330757803522,45.86836716450355,7.344005520526261,30.62325117856182,10.1501531372438,15.51248839313582,17.97404495263578,22.84164090774314),
    @(14.58704262209039,11.34696419554614,6.332928971172755,12.06876259344203,45.82674432728718,7.344005520526261,30.6230105739433,10.15532687947678,15.45089519528917,17.94811253876559,22.86121931396401),
    @(14.24897966337062,11.11325229404649,6.310013836451072,12.05484408481507,45.68486307672509,7.344005520526261,30.63122264709624,10.1768516238737,15.2215198196991,17.85548987395009,22.9387140739038),
    @(13.88850200416772,10.86659881047817,6.288771260488205,12.04692558802496,45.56581523101181,7.344005520526261,30.65955843543717,10.20420087244857,14.98357797878284,17.76783943840659,23.0294617783634)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
